$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.392.26"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.842.96"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.89"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6310"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07524"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2924"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.40"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D12").Value = "1.843.48"
$ws.Range("E12").Value = "  -7.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.999"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6790"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001038"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "2.092.36"
$ws.Range("E17").Value = "  -7.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.165"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "29.423.34"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.11"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.443"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.01"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.365"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.57"
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.458"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.099"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.028"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.841"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.157"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7120"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.589"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "1.245.40"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.771"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.343"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9016"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.68"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.72"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.114"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000118"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3996"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.948"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.671"
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1122"
$ws.Range("E51").Value = "  -0.29%  "
